$d = $word.ActiveDocument

$d.Content.Find.Execute("Total unique exceptions: 1591", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Total unique exceptions: 662", 2)

$d.Content.Find.Execute("Exceptions in BOTH: 825", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Exceptions in BOTH: 243", 2)

$d.Content.Find.Execute("Exceptions OXFORD only: 395", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Exceptions OXFORD only: 282", 2)

$d.Content.Find.Execute("Exceptions RIVERSIDE only: 371", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Exceptions RIVERSIDE only: 137", 2)
